# Updated symbol list on Wed Dec 14 13:57:23 UTC 2022 with GitHub Actions
#
# Column D on Sheet1 holds crypto "Price" values. In the source workbook
# they are stored as plain text (not numbers) even though they look
# numeric (e.g. "269.63"). Assigning a numeric-looking string straight to
# Range.Value would make Excel auto-convert it to a real number, which
# would not match the original text formatting. To keep the values as
# text we momentarily force the cell to Text number format ("@"), write
# the new value, then clear the formatting back to the default (General)
# so no visible formatting change is left behind - matching how the
# workbook looked before the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextPrice([string]$cellRef, [string]$newValue) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $newValue
    $cell.ClearFormats()
}

Set-TextPrice "D2"  "269.66"
Set-TextPrice "D3"  "22.91"
Set-TextPrice "D4"  "6.369"
Set-TextPrice "D6"  "3.640"
Set-TextPrice "D7"  "6.696"
Set-TextPrice "D8"  "1.374"
Set-TextPrice "D9"  "0.8342"
Set-TextPrice "D11" "0.1631"
Set-TextPrice "D12" "0.08390"
Set-TextPrice "D13" "0.03417"
Set-TextPrice "D14" "0.03111"
Set-TextPrice "D15" "0.09318"
Set-TextPrice "D16" "3.890"
Set-TextPrice "D17" "0.001710"
Set-TextPrice "D18" "0.04810"
Set-TextPrice "D19" "0.006265"
Set-TextPrice "D20" "0.001088"
Set-TextPrice "D21" "0.003328"
Set-TextPrice "D22" "0.0001500"
Set-TextPrice "D23" "3.737"
Set-TextPrice "D24" "2.371"
Set-TextPrice "D40" "0.04682"
Set-TextPrice "D41" "0.006901"
Set-TextPrice "D43" "0.003340"
Set-TextPrice "D44" "0.01118"
Set-TextPrice "D45" "0.00006256"
Set-TextPrice "D46" "0.00000000750"
Set-TextPrice "D47" "0.9002"
Set-TextPrice "D48" "0.07471"
Set-TextPrice "D49" "0.00001400"
Set-TextPrice "D50" "0.01240"
